$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 29.253501
$ws.Range("H2").Value = 87.760503
$ws.Range("I2").Value = 0.7876335333413836
$ws.Range("J2").Value = 0.7876335333413838
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 23.63579766666667
$ws.Range("N2").Value = 70.907393
$ws.Range("O2").Value = 0.06827844587621175
$ws.Range("P2").Value = 0.06827844587621175
$ws.Range("Q2").Value = 691.429830677631
$ws.Range("R2").Value = 6222.868476098679
$ws.Range("S2").Value = 0.05377839357653909
$ws.Range("T2").Value = 0.05377839357653909
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 29.253501
$ws.Range("H3").Value = 87.760503
$ws.Range("I3").Value = 0.7876335333413836
$ws.Range("J3").Value = 0.7876335333413838
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 181.2883913333334
$ws.Range("N3").Value = 543.865174
$ws.Range("O3").Value = 0.5237009467675041
$ws.Range("P3").Value = 0.523700946767504
$ws.Range("Q3").Value = 5303.320137158058
$ws.Range("R3").Value = 47729.88123442252
$ws.Range("S3").Value = 0.4124844271167171
$ws.Range("T3").Value = 0.4124844271167171
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 29.253501
$ws.Range("H4").Value = 87.760503
$ws.Range("I4").Value = 0.7876335333413836
$ws.Range("J4").Value = 0.7876335333413838
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 111.1005463333333
$ws.Range("N4").Value = 333.301639
$ws.Range("O4").Value = 0.3209442197221123
$ws.Range("P4").Value = 0.3209442197221123
$ws.Range("Q4").Value = 3250.079943262713
$ws.Range("R4").Value = 29250.71948936442
$ws.Range("S4").Value = 0.2527864297852207
$ws.Range("T4").Value = 0.2527864297852207
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 29.253501
$ws.Range("H5").Value = 87.760503
$ws.Range("I5").Value = 0.7876335333413836
$ws.Range("J5").Value = 0.7876335333413838
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 30.14303933333333
$ws.Range("N5").Value = 90.42911799999999
$ws.Range("O5").Value = 0.08707638763417187
$ws.Range("P5").Value = 0.08707638763417187
$ws.Range("Q5").Value = 881.7894312807059
$ws.Range("R5").Value = 7936.104881526353
$ws.Range("S5").Value = 0.06858428286290676
$ws.Range("T5").Value = 0.06858428286290677
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 4.749137666666667
$ws.Range("H6").Value = 14.247413
$ws.Range("I6").Value = 0.1278677748937237
$ws.Range("J6").Value = 0.1278677748937237
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 23.63579766666667
$ws.Range("N6").Value = 70.907393
$ws.Range("O6").Value = 0.06827844587621175
$ws.Range("P6").Value = 0.06827844587621175
$ws.Range("Q6").Value = 112.2496569804788
$ws.Range("R6").Value = 1010.246912824309
$ws.Range("S6").Value = 0.008730612947392741
$ws.Range("T6").Value = 0.008730612947392743
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 4.749137666666667
$ws.Range("H7").Value = 14.247413
$ws.Range("I7").Value = 0.1278677748937237
$ws.Range("J7").Value = 0.1278677748937237
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 181.2883913333334
$ws.Range("N7").Value = 543.865174
$ws.Range("O7").Value = 0.5237009467675041
$ws.Range("P7").Value = 0.523700946767504
$ws.Range("Q7").Value = 860.9635278105403
$ws.Range("R7").Value = 7748.671750294863
$ws.Range("S7").Value = 0.06696447477289719
$ws.Range("T7").Value = 0.06696447477289719
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 4.749137666666667
$ws.Range("H8").Value = 14.247413
$ws.Range("I8").Value = 0.1278677748937237
$ws.Range("J8").Value = 0.1278677748937237
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 111.1005463333333
$ws.Range("N8").Value = 333.301639
$ws.Range("O8").Value = 0.3209442197221123
$ws.Range("P8").Value = 0.3209442197221123
$ws.Range("Q8").Value = 527.6317893788786
$ws.Range("R8").Value = 4748.686104409907
$ws.Range("S8").Value = 0.04103842324086885
$ws.Range("T8").Value = 0.04103842324086886
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 4.749137666666667
$ws.Range("H9").Value = 14.247413
$ws.Range("I9").Value = 0.1278677748937237
$ws.Range("J9").Value = 0.1278677748937237
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 30.14303933333333
$ws.Range("N9").Value = 90.42911799999999
$ws.Range("O9").Value = 0.08707638763417187
$ws.Range("P9").Value = 0.08707638763417187
$ws.Range("Q9").Value = 143.1534434857482
$ws.Range("R9").Value = 1288.380991371734
$ws.Range("S9").Value = 0.01113426393256491
$ws.Range("T9").Value = 0.01113426393256492
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 3.138366333333334
$ws.Range("H10").Value = 9.415099000000001
$ws.Range("I10").Value = 0.08449869176489255
$ws.Range("J10").Value = 0.08449869176489258
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 23.63579766666667
$ws.Range("N10").Value = 70.907393
$ws.Range("O10").Value = 0.06827844587621175
$ws.Range("P10").Value = 0.06827844587621175
$ws.Range("Q10").Value = 74.17779165854523
$ws.Range("R10").Value = 667.6001249269071
$ws.Range("S10").Value = 0.005769439352279916
$ws.Range("T10").Value = 0.005769439352279917
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 3.138366333333334
$ws.Range("H11").Value = 9.415099000000001
$ws.Range("I11").Value = 0.08449869176489255
$ws.Range("J11").Value = 0.08449869176489258
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 181.2883913333334
$ws.Range("N11").Value = 543.865174
$ws.Range("O11").Value = 0.5237009467675041
$ws.Range("P11").Value = 0.523700946767504
$ws.Range("Q11").Value = 568.949383984692
$ws.Range("R11").Value = 5120.544455862227
$ws.Range("S11").Value = 0.04425204487788973
$ws.Range("T11").Value = 0.04425204487788974
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 3.138366333333334
$ws.Range("H12").Value = 9.415099000000001
$ws.Range("I12").Value = 0.08449869176489255
$ws.Range("J12").Value = 0.08449869176489258
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 111.1005463333333
$ws.Range("N12").Value = 333.301639
$ws.Range("O12").Value = 0.3209442197221123
$ws.Range("P12").Value = 0.3209442197221123
$ws.Range("Q12").Value = 348.6742142274735
$ws.Range("R12").Value = 3138.067928047261
$ws.Range("S12").Value = 0.02711936669602271
$ws.Range("T12").Value = 0.02711936669602272
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 3.138366333333334
$ws.Range("H13").Value = 9.415099000000001
$ws.Range("I13").Value = 0.08449869176489255
$ws.Range("J13").Value = 0.08449869176489258
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 30.14303933333333
$ws.Range("N13").Value = 90.42911799999999
$ws.Range("O13").Value = 0.08707638763417187
$ws.Range("P13").Value = 0.08707638763417187
$ws.Range("Q13").Value = 94.59989982807578
$ws.Range("R13").Value = 851.399098452682
$ws.Range("S13").Value = 0.007357840838700191
$ws.Range("T13").Value = 0.007357840838700193
